$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.171443343162537
$ws.Range("B1").Value = 2.22373366355896
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.330950021743774
$ws.Range("E1").Value = 1.22807776927948
